# C1--C2-and-C3-PowerPoint.pptx edit:
#   1) Re-style the cash-flow table on slide 16 with a different built-in
#      PowerPoint table style (tableStyleId GUID change).
#   2) Re-colour the deck's theme from the imported "Integral" palette to
#      the stock "Office Theme" palette (the font scheme / format scheme
#      were already identical between the two themes, so only the 12
#      scheme colours need to move).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 -------------------------------------------
$slide16 = $p.Slides.Item(16)
$tableShape = $slide16.Shapes.Item(3)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{A61E07E4-60A8-4FE8-9FE6-A7B9999DE038}", $false)
}

# --- 2) Theme colours -------------------------------------------------------
# Office Theme scheme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM RGB() integers (0x00BBGGRR).
$themeColors = $p.Slides.Item(1).ThemeColorScheme
$themeColors.Colors(1).RGB  = 0         # dk1      000000
$themeColors.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$themeColors.Colors(3).RGB  = 6968388   # dk2      44546A
$themeColors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$themeColors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$themeColors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$themeColors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$themeColors.Colors(8).RGB  = 49407     # accent4  FFC000
$themeColors.Colors(9).RGB  = 12874308  # accent5  4472C4
$themeColors.Colors(10).RGB = 4697456   # accent6  70AD47
$themeColors.Colors(11).RGB = 12673797  # hlink    0563C1
$themeColors.Colors(12).RGB = 7491477   # folHlink 954F72
